# Apply updated dSF (column F) values as part of a data repull / push.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new dSF (column F) value
$updates = @{
    3  = 1
    4  = 5
    5  = 2
    6  = 1
    7  = 1
    8  = -1
    9  = 3
    10 = 4
    12 = 2
    13 = -3
    16 = -8
    17 = 1
    18 = 4
    20 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
